$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 and Row 6 contain two separate observation records that need to
# swap places (all field values exchanged), while staying on the same
# physical spreadsheet rows (formatting/row number unchanged).

$cols = @("A","B","D","E","F","G","H","P","Q","R","Z","AB")

foreach ($col in $cols) {
    $addr5 = "$col" + "5"
    $addr6 = "$col" + "6"
    $val5 = $ws.Range($addr5).Value2
    $val6 = $ws.Range($addr6).Value2
    $ws.Range($addr5).Value = $val6
    $ws.Range($addr6).Value = $val5
}
